# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Update existing F-column (quantity) values for rows 2-4, fix row 5 to be
# DC_001 (instead of DC_002) with an updated quantity/horizon, and append two
# new rows (6 and 7) for MAT_B at DC_002 and PLANT_001.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NetDemand")

# --- Update quantity values for existing rows (MAT_A) ---
$ws.Range("F2").Value = -225
$ws.Range("F3").Value = -422
$ws.Range("F4").Value = -709

# --- Row 5 (MAT_B): location DC_002 -> DC_001, quantity + horizon_days updated ---
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -98
$ws.Range("H5").Value = 4

# --- New row 6 (MAT_B, DC_002) ---
$ws.Range("A6").Value = "MAT_B"
$ws.Range("B6").Value = "DC_002"
$ws.Range("C6").Value = $ws.Range("C5").Value2
$ws.Range("C6").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("D6").Value = "Distribution Demand - Forecast"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -48
$ws.Range("G6").Value = $ws.Range("G5").Value2
$ws.Range("G6").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("H6").Value = 1

# --- New row 7 (MAT_B, PLANT_001) ---
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = $ws.Range("C5").Value2
$ws.Range("C7").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -98
$ws.Range("G7").Value = $ws.Range("G5").Value2
$ws.Range("G7").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("H7").Value = 1
